# Updates cached market-price / profit figures on each job sheet.
# Generated from the scheduled-runner diff (H,I,J,K,L,M,N columns).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")

$ws_ALC.Range("H88").Value = 3175.1538
$ws_ALC.Range("I88").Value = 3898.6667
$ws_ALC.Range("J88").Value = 2958.1
$ws_ALC.Range("K88").Value = 3898.6667
$ws_ALC.Range("L88").Value = 2958.1
$ws_ALC.Range("M88").Value = -3492.6667
$ws_ALC.Range("N88").Value = -3770.1

$ws_ALC.Range("H91").Value = 3175.1538
$ws_ALC.Range("I91").Value = 3898.6667
$ws_ALC.Range("J91").Value = 2958.1
$ws_ALC.Range("K91").Value = 3898.6667
$ws_ALC.Range("L91").Value = 2958.1
$ws_ALC.Range("M91").Value = -2494.6667
$ws_ALC.Range("N91").Value = -5766.1

$ws_ALC.Range("H116").Value = 14916.632
$ws_ALC.Range("I116").Value = 17280.8
$ws_ALC.Range("K116").Value = 17280.8
$ws_ALC.Range("M116").Value = -13838.8

$ws_ALC.Range("H135").Value = 2730.5557
$ws_ALC.Range("I135").Value = 1596.3334
$ws_ALC.Range("K135").Value = 14367.0006
$ws_ALC.Range("M135").Value = -11832.0006

$ws_ALC.Range("H137").Value = 1737947.8
$ws_ALC.Range("I137").Value = 1791.6875
$ws_ALC.Range("J137").Value = 3474103.8
$ws_ALC.Range("K137").Value = 5375.0625
$ws_ALC.Range("L137").Value = 10422311.4
$ws_ALC.Range("M137").Value = -2825.0625
$ws_ALC.Range("N137").Value = -10427411.4

$ws_ALC.Range("H138").Value = 2374.78
$ws_ALC.Range("I138").Value = 1163.381
$ws_ALC.Range("J138").Value = 2696.7974
$ws_ALC.Range("K138").Value = 3490.143
$ws_ALC.Range("L138").Value = 8090.3922
$ws_ALC.Range("M138").Value = 1649.857
$ws_ALC.Range("N138").Value = -18370.3922

$ws_ARM = $wb.Worksheets.Item("ARM")

$ws_ARM.Range("H110").Value = 1819
$ws_ARM.Range("I110").Value = 1795.2222
$ws_ARM.Range("J110").Value = 1961.6666
$ws_ARM.Range("K110").Value = 1795.2222
$ws_ARM.Range("L110").Value = 1961.6666
$ws_ARM.Range("M110").Value = 249.7778000000001
$ws_ARM.Range("N110").Value = -6051.6666

$ws_BSM = $wb.Worksheets.Item("BSM")

$ws_BSM.Range("H134").Value = 2307107.5
$ws_BSM.Range("I134").Value = 2859176.8
$ws_BSM.Range("K134").Value = 8577530.399999999
$ws_BSM.Range("M134").Value = -8574995.399999999

$ws_CRP = $wb.Worksheets.Item("CRP")

$ws_CRP.Range("H7").Value = 35.411766
$ws_CRP.Range("I7").Value = 31.4375
$ws_CRP.Range("J7").Value = 99
$ws_CRP.Range("K7").Value = 31.4375
$ws_CRP.Range("L7").Value = 99
$ws_CRP.Range("M7").Value = 81.5625
$ws_CRP.Range("N7").Value = -325

$ws_CRP.Range("H22").Value = 250.06667
$ws_CRP.Range("I22").Value = 232.21428
$ws_CRP.Range("K22").Value = 232.21428
$ws_CRP.Range("M22").Value = 117.78572

$ws_CRP.Range("H31").Value = 4626.2856
$ws_CRP.Range("J31").Value = 6549.45
$ws_CRP.Range("L31").Value = 6549.45
$ws_CRP.Range("N31").Value = -7139.45

$ws_CRP.Range("H34").Value = 4626.2856
$ws_CRP.Range("J34").Value = 6549.45
$ws_CRP.Range("L34").Value = 6549.45
$ws_CRP.Range("N34").Value = -6953.45

$ws_CRP.Range("H58").Value = 2348.6606
$ws_CRP.Range("I58").Value = 2142.848
$ws_CRP.Range("J58").Value = 3295.4
$ws_CRP.Range("K58").Value = 2142.848
$ws_CRP.Range("L58").Value = 3295.4
$ws_CRP.Range("M58").Value = -1939.848
$ws_CRP.Range("N58").Value = -3701.4

$ws_CRP.Range("H122").Value = 4392.3687
$ws_CRP.Range("I122").Value = 4171.1816
$ws_CRP.Range("J122").Value = 4696.5
$ws_CRP.Range("K122").Value = 12513.5448
$ws_CRP.Range("L122").Value = 14089.5
$ws_CRP.Range("M122").Value = -10063.5448
$ws_CRP.Range("N122").Value = -18989.5

$ws_CRP.Range("H134").Value = 3130.0557
$ws_CRP.Range("I134").Value = 3021.3125
$ws_CRP.Range("K134").Value = 9063.9375
$ws_CRP.Range("M134").Value = -6528.9375

$ws_CRP.Range("H136").Value = 2348.6606
$ws_CRP.Range("I136").Value = 2142.848
$ws_CRP.Range("J136").Value = 3295.4
$ws_CRP.Range("K136").Value = 6428.544
$ws_CRP.Range("L136").Value = 9886.200000000001
$ws_CRP.Range("M136").Value = -3878.544
$ws_CRP.Range("N136").Value = -14986.2

$ws_CUL = $wb.Worksheets.Item("CUL")

$ws_CUL.Range("H5").Value = 1921.9231
$ws_CUL.Range("I5").Value = 1221.2222
$ws_CUL.Range("J5").Value = 3498.5
$ws_CUL.Range("K5").Value = 3663.6666
$ws_CUL.Range("L5").Value = 10495.5
$ws_CUL.Range("M5").Value = -3551.6666
$ws_CUL.Range("N5").Value = -10719.5

$ws_CUL.Range("H18").Value = 4287.8
$ws_CUL.Range("I18").Value = 3146.5
$ws_CUL.Range("K18").Value = 9439.5
$ws_CUL.Range("M18").Value = -9270.5

$ws_CUL.Range("H64").Value = 100
$ws_CUL.Range("I64").Value = 100
$ws_CUL.Range("K64").Value = 300
$ws_CUL.Range("M64").Value = -30

$ws_CUL.Range("H67").Value = 100
$ws_CUL.Range("I67").Value = 100
$ws_CUL.Range("K67").Value = 300
$ws_CUL.Range("M67").Value = 636

$ws_CUL.Range("H113").Value = 1448.5714
$ws_CUL.Range("I113").Value = 656
$ws_CUL.Range("K113").Value = 1968
$ws_CUL.Range("M113").Value = 202

$ws_CUL.Range("H117").Value = 3782.2666
$ws_CUL.Range("I117").Value = 2846
$ws_CUL.Range("J117").Value = 4122.727
$ws_CUL.Range("K117").Value = 8538
$ws_CUL.Range("L117").Value = 12368.181
$ws_CUL.Range("M117").Value = -5096
$ws_CUL.Range("N117").Value = -19252.181

$ws_CUL.Range("H120").Value = 20642.309
$ws_CUL.Range("I120").Value = 10837.5
$ws_CUL.Range("J120").Value = 25000
$ws_CUL.Range("K120").Value = 32512.5
$ws_CUL.Range("L120").Value = 75000
$ws_CUL.Range("M120").Value = -27674.5
$ws_CUL.Range("N120").Value = -84676

$ws_CUL.Range("H135").Value = 1921.9231
$ws_CUL.Range("I135").Value = 1221.2222
$ws_CUL.Range("J135").Value = 3498.5
$ws_CUL.Range("K135").Value = 10990.9998
$ws_CUL.Range("L135").Value = 31486.5
$ws_CUL.Range("M135").Value = -8455.9998
$ws_CUL.Range("N135").Value = -36556.5

$ws_GSM = $wb.Worksheets.Item("GSM")

$ws_GSM.Range("H52").Value = 46025.285
$ws_GSM.Range("J52").Value = 45470.75
$ws_GSM.Range("L52").Value = 45470.75
$ws_GSM.Range("N52").Value = -45988.75

$ws_GSM.Range("H102").Value = 2784.0715
$ws_GSM.Range("I102").Value = 2607
$ws_GSM.Range("J102").Value = 3433.3333
$ws_GSM.Range("K102").Value = 2607
$ws_GSM.Range("L102").Value = 3433.3333
$ws_GSM.Range("M102").Value = -985
$ws_GSM.Range("N102").Value = -6677.3333

$ws_LTW = $wb.Worksheets.Item("LTW")

$ws_LTW.Range("H7").Value = 13238.571
$ws_LTW.Range("I7").Value = 14278.333
$ws_LTW.Range("K7").Value = 14278.333
$ws_LTW.Range("M7").Value = -14166.333

$ws_LTW.Range("H16").Value = 284.85715
$ws_LTW.Range("I16").Value = 284.85715
$ws_LTW.Range("J16").Value = 0
$ws_LTW.Range("K16").Value = 284.85715
$ws_LTW.Range("L16").Value = 0
$ws_LTW.Range("N16").Value = -114.85715
$ws_LTW.Range("M16").Value = ""

$ws_LTW.Range("H40").Value = 4086.7693
$ws_LTW.Range("J40").Value = 6200
$ws_LTW.Range("L40").Value = 6200
$ws_LTW.Range("N40").Value = -6472

$ws_LTW.Range("H122").Value = 19187.6
$ws_LTW.Range("I122").Value = 21523.2
$ws_LTW.Range("J122").Value = 12180.8
$ws_LTW.Range("K122").Value = 64569.60000000001
$ws_LTW.Range("L122").Value = 36542.39999999999
$ws_LTW.Range("M122").Value = -62119.60000000001
$ws_LTW.Range("N122").Value = -41442.39999999999

$ws_LTW.Range("H126").Value = 13238.571
$ws_LTW.Range("I126").Value = 14278.333
$ws_LTW.Range("K126").Value = 42834.999
$ws_LTW.Range("M126").Value = -40364.999

$ws_LTW.Range("H132").Value = 772895.3
$ws_LTW.Range("I132").Value = 1114026.9
$ws_LTW.Range("J132").Value = 5349.25
$ws_LTW.Range("K132").Value = 3342080.7
$ws_LTW.Range("L132").Value = 16047.75
$ws_LTW.Range("M132").Value = -3339550.7
$ws_LTW.Range("N132").Value = -21107.75

$ws_LTW.Range("H136").Value = 24999.8
$ws_LTW.Range("J136").Value = 24999.8
$ws_LTW.Range("L136").Value = 74999.39999999999
$ws_LTW.Range("N136").Value = -80099.39999999999

$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_WVR.Range("H96").Value = 9624.6
$ws_WVR.Range("I96").Value = 4401.8
$ws_WVR.Range("J96").Value = 12236
$ws_WVR.Range("K96").Value = 4401.8
$ws_WVR.Range("L96").Value = 12236
$ws_WVR.Range("M96").Value = -3028.8
$ws_WVR.Range("N96").Value = -14982

$ws_WVR.Range("H107").Value = 336.875
$ws_WVR.Range("I107").Value = 273.57144
$ws_WVR.Range("K107").Value = 820.71432
$ws_WVR.Range("M107").Value = 1099.28568

$ws_WVR.Range("H136").Value = 2355.4
$ws_WVR.Range("I136").Value = 1712
$ws_WVR.Range("J136").Value = 4124.75
$ws_WVR.Range("K136").Value = 5136
$ws_WVR.Range("L136").Value = 12374.25
$ws_WVR.Range("M136").Value = -2586
$ws_WVR.Range("N136").Value = -17474.25
